# "Generate Report for Handback"
#
# Refresh the handback-status report: update the "Latest HO Xliff Generate
# Date" / "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the file that was just handed back
# (ef514347-a333-4a59-8b2c-bca73fe2522e.md), on the Overview sheet as well
# as on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-02 22:54:10"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-02 22:53:59"
$zhcn.Range("K2").Value = "2016-09-02 22:54:30"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-02 22:54:10"
$dede.Range("K2").Value = "2016-09-02 22:54:37"
